$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E28:F29").NumberFormat = "@"

$ws.Range("A28").Value = "2024-09-12 19:36:57"
$ws.Range("B28").Value = "check_availability"
$ws.Range("C28").Value = "MOCKURL_https://www.opentable.com/r/bar-spero-washington/"
$ws.Range("D28").Value = "MOCK_No availability for the selected date."
$ws.Range("E28").Value = "2024-09-12"
$ws.Range("F28").Value = "19:36:57"

$ws.Range("A29").Value = "2024-09-12 19:42:53"
$ws.Range("B29").Value = "check_availability"
$ws.Range("C29").Value = "MOCKURL_https://www.opentable.com/r/bar-spero-washington/"
$ws.Range("D29").Value = "MOCK_No availability for the selected date."
$ws.Range("E29").Value = "2024-09-12"
$ws.Range("F29").Value = "19:42:53"
